$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Update "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) for rows 3 and 4
# (new progress as of 04-Nov-2025).
#
# The "LAST UPDATE" column holds its date as literal text (e.g. "03-Nov-2025"),
# not a real date value, so we force a Text number format before writing the
# new value - otherwise Excel would auto-convert the date-like string into a
# date serial number. We restore the General format afterwards to keep the
# cell formatting as close as possible to how it started.
$ws.Range("H3").Value = 521
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "04-Nov-2025"
$ws.Range("I3").NumberFormat = "General"

$ws.Range("H4").Value = 521
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "04-Nov-2025"
$ws.Range("I4").NumberFormat = "General"
